$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-unused new_twigs/old_twigs columns (G:H) for header + data rows
$ws.Range("G1:H11").ClearContents()

# Update new_leaves / old_leaves measured data (columns E/F) for rows 2-11
$ws.Range("E2").Value = 0.06421731
$ws.Range("F2").Value = 0.10140037
$ws.Range("E3").Value = 0.07920631
$ws.Range("F3").Value = 0.11555194
$ws.Range("E4").Value = 0.06587365
$ws.Range("F4").Value = 0.10388062
$ws.Range("E5").Value = 0.07000135
$ws.Range("F5").Value = 0.11193648
$ws.Range("E6").Value = 0.07652952
$ws.Range("F6").Value = 0.12051627
$ws.Range("E7").Value = 0.07245142
$ws.Range("F7").Value = 0.11830346
$ws.Range("E8").Value = 0.0632167
$ws.Range("F8").Value = 0.10009229
$ws.Range("E9").Value = 0.04721238
$ws.Range("F9").Value = 0.0647876
$ws.Range("E10").Value = 0.05082157
$ws.Range("F10").Value = 0.07828481
$ws.Range("E11").Value = 0.06810342
$ws.Range("F11").Value = 0.10806316

# Update New_max (B16) and old_max (B18) -- this ripples into C2:D11 formulas
$ws.Range("B16").Value = 0.02909091
$ws.Range("B18").Value = 0.08641584

# Remove the stray helper cell in row 26
$ws.Range("C26").ClearContents()

# Add the new note about the correction, in a fresh row below the existing note
$ws.Range("A23").Value = "Sto a scoppia', ho cambiato di ogni, non so come cazzo avevo ragionato prima, ma il ragionamento corretto dovrebbe essere: qui mi sono ricalcolato le masse dei rametti partendo dalle pesate/ lunghezze di laboratorio e rapportandoli ai rami trattati (ma in teorira, spero, solo per proporzionarli tra loro all'interno dei valori da laboratorio. Puo' essere che ho fatto una cazzata ad eliminare i lavori di lab piu alti che in effetti dovrebbero rappresentare i rami piu grossi, ma li avevo tolti quando avevo calcolato erroneamente le masse dei rametti con un doppio passaggio concettualmente sbagliato"

# Update the active selection to match the edited range
$ws.Range("D2:D11").Select()
